$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay as literal text (they are
# free-form "28.499.95" style strings, not real numbers), matching
# the source data which stores them as inline/shared strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.499.95'
$ws.Range("E2").Value = '  +0.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.50'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.86'
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5165'
$ws.Range("E7").Value = '  +2.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3869'
$ws.Range("E8").Value = '  -1.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08326'
$ws.Range("E9").Value = '  +8.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.01'
$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.120'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.409'
$ws.Range("E12").Value = '  +2.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.16'
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.498'
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.825.62'
$ws.Range("E16").Value = '  +0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.98'
$ws.Range("E17").Value = '  +0.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001125'
$ws.Range("E18").Value = '  +4.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06646'
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.76'
$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.043'
$ws.Range("E22").Value = '  -1.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.546.76'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  +2.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.298'
$ws.Range("E25").Value = '  +1.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.16'
$ws.Range("E26").Value = '  +2.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.87'
$ws.Range("E27").Value = '  +1.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.033.89'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.400'
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.98'
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1092'
$ws.Range("E31").Value = '  +0.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.094'
$ws.Range("E32").Value = '  -3.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.07608'
$ws.Range("E33").Value = '  +7.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.731'
$ws.Range("E34").Value = '  +1.41%  '

$ws.Range("E35").Value = '  +0.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2228'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("E37").Value = '  +2.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.281'
$ws.Range("E38").Value = '  +3.07%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.75'
$ws.Range("E39").Value = '  +5.00%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.738'
$ws.Range("E40").Value = '  -2.17%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6383'
$ws.Range("E41").Value = '  +2.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.193'
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.398'
$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.51'
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6110'
$ws.Range("E45").Value = '  +3.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.800'
$ws.Range("E46").Value = '  +2.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.67'
$ws.Range("E47").Value = '  +2.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.995'
$ws.Range("E48").Value = '  +1.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.206'
$ws.Range("E49").Value = '  +1.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06996'
$ws.Range("E50").Value = '  +0.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.28'
$ws.Range("E51").Value = '  +0.82%  '
